$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Week 7 topic cell: "Testing and Rework, MidTerm" -> "Testing and Rework"
#    (the "MidTerm" portion moves down to the Week 9 topic cell below)
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("Testing and Rework, MidTerm", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Testing and Rework", 2)

# ------------------------------------------------------------------
# 2. Remove the "20% midterm" paragraph from the Week 7 Assessment/Weight
#    cell (it currently sits above "15% Assign. 2").
# ------------------------------------------------------------------
$midtermPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("20% midterm")) {
        $midtermPara = $p
        break
    }
}
if ($midtermPara -ne $null) {
    $midtermPara.Range.Delete()
}

# ------------------------------------------------------------------
# 3. Insert a new "20% midterm" paragraph in the Week 8 Assessment/Weight
#    cell, right before the existing "(5% Assign. 3 part B)" paragraph.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$weekFoundRow = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    try {
        $cellText = $t.Cell($r, 5).Range.Text
    } catch {
        continue
    }
    if ($cellText.StartsWith("(5% Assign. 3 part B)")) {
        $weekFoundRow = $r
        break
    }
}
$assignBCell = $t.Cell($weekFoundRow, 5)
$assignBCell.Range.InsertBefore("20% midterm`r")

# ------------------------------------------------------------------
# 4. Week 9 topic cell: "Teamwork" -> "Teamwork," followed by a new
#    centered paragraph containing "MidTerm".
# ------------------------------------------------------------------
$t2 = $d.Tables.Item(1)
$teamworkRow = -1
for ($r = 1; $r -le $t2.Rows.Count; $r++) {
    try {
        $cellText = $t2.Cell($r, 2).Range.Text
    } catch {
        continue
    }
    if ($cellText.StartsWith("Teamwork")) {
        $teamworkRow = $r
        break
    }
}
$teamworkCell = $t2.Cell($teamworkRow, 2)
$teamworkCell.Range.InsertAfter(",`rMidTerm")

Write-Output "edits applied"
